# Partner import template: mark the mandatory header columns with a
# trailing "(*)" suffix, matching the updated import/catalogue rules.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "Name Abr(*)"
$ws.Range("B1").Value = "Name EN(*)"
$ws.Range("C1").Value = "Name Local(*)"
$ws.Range("D1").Value = "Taxcode(*)"
$ws.Range("E1").Value = "Category(*)"
$ws.Range("H1").Value = "Billing  Address EN(*)"
$ws.Range("I1").Value = "Billing Address Local(*)"
$ws.Range("J1").Value = " City (*)"
$ws.Range("K1").Value = "Country(*)"
$ws.Range("M1").Value = "Shipping  Address EN(*)"
$ws.Range("N1").Value = "Shipping Address Local(*)"
$ws.Range("O1").Value = " City (*)"
$ws.Range("P1").Value = "Country(*)"

# Move the active selection the way the author last left the sheet.
$ws.Range("U15").Select()
